$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 550
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 550
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 550
$ws.Range("N12").Value = -890
$ws.Range("H38").Value = 302.5
$ws.Range("I38").Value = 302.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 907.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -535.5
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 2026.625
$ws.Range("I58").Value = 57.5
$ws.Range("J58").Value = 2683
$ws.Range("K58").Value = 172.5
$ws.Range("L58").Value = 8049
$ws.Range("M58").Value = -22.5
$ws.Range("N58").Value = -8349
$ws.Range("H80").Value = 346.27274
$ws.Range("I80").Value = 285
$ws.Range("J80").Value = 381.2857
$ws.Range("K80").Value = 855
$ws.Range("L80").Value = 1143.8571
$ws.Range("M80").Value = 143
$ws.Range("N80").Value = -3139.8571
$ws.Range("H83").Value = 346.27274
$ws.Range("I83").Value = 285
$ws.Range("J83").Value = 381.2857
$ws.Range("K83").Value = 2565
$ws.Range("L83").Value = 3431.5713
$ws.Range("M83").Value = 2427
$ws.Range("N83").Value = -13415.5713
$ws.Range("H86").Value = 1382
$ws.Range("I86").Value = 949.6667
$ws.Range("J86").Value = 1900.8
$ws.Range("K86").Value = 949.6667
$ws.Range("L86").Value = 1900.8
$ws.Range("M86").Value = 173.3333
$ws.Range("N86").Value = -4146.8
$ws.Range("H89").Value = 1382
$ws.Range("I89").Value = 949.6667
$ws.Range("J89").Value = 1900.8
$ws.Range("K89").Value = 4748.3335
$ws.Range("L89").Value = 9504
$ws.Range("M89").Value = 867.6665000000003
$ws.Range("N89").Value = -20736
$ws.Range("H100").Value = 5500
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 5500
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 5500
$ws.Range("N100").Value = -6582
$ws.Range("H132").Value = 11362.65
$ws.Range("I132").Value = 11362.65
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 34087.95
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -31557.95
$ws.Range("H138").Value = 8556.429
$ws.Range("I138").Value = 348.5
$ws.Range("J138").Value = 11839.6
$ws.Range("K138").Value = 1045.5
$ws.Range("L138").Value = 35518.8
$ws.Range("M138").Value = 4094.5
$ws.Range("N138").Value = -45798.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2500
$ws.Range("I2").Value = 2500
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2387
$ws.Range("H8").Value = 4233
$ws.Range("I8").Value = 150
$ws.Range("J8").Value = 5866.2
$ws.Range("K8").Value = 150
$ws.Range("L8").Value = 5866.2
$ws.Range("M8").Value = -6
$ws.Range("N8").Value = -6154.2
$ws.Range("H11").Value = 8002129.5
$ws.Range("I11").Value = 3668.6667
$ws.Range("J11").Value = 14000975
$ws.Range("K11").Value = 3668.6667
$ws.Range("L11").Value = 14000975
$ws.Range("M11").Value = -3524.6667
$ws.Range("N11").Value = -14001263
$ws.Range("H13").Value = 1950
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1950
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1950
$ws.Range("N13").Value = -2238
$ws.Range("M13").ClearContents()
$ws.Range("H45").Value = 3216.25
$ws.Range("I45").Value = 2177.2
$ws.Range("J45").Value = 4948
$ws.Range("K45").Value = 2177.2
$ws.Range("L45").Value = 4948
$ws.Range("M45").Value = -1800.2
$ws.Range("N45").Value = -5702
$ws.Range("H88").Value = 1453.8182
$ws.Range("I88").Value = 1368.7142
$ws.Range("J88").Value = 1602.75
$ws.Range("K88").Value = 1368.7142
$ws.Range("L88").Value = 1602.75
$ws.Range("M88").Value = -962.7141999999999
$ws.Range("N88").Value = -2414.75
$ws.Range("H91").Value = 1453.8182
$ws.Range("I91").Value = 1368.7142
$ws.Range("J91").Value = 1602.75
$ws.Range("K91").Value = 1368.7142
$ws.Range("L91").Value = 1602.75
$ws.Range("M91").Value = 35.28580000000011
$ws.Range("N91").Value = -4410.75
$ws.Range("H97").Value = 1080.25
$ws.Range("I97").Value = 1112.4445
$ws.Range("J97").Value = 983.6667
$ws.Range("K97").Value = 1112.4445
$ws.Range("L97").Value = 983.6667
$ws.Range("M97").Value = -616.4445000000001
$ws.Range("N97").Value = -1975.6667
$ws.Range("H116").Value = 2500
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -206
$ws.Range("H132").Value = 1186.6562
$ws.Range("I132").Value = 1152.7037
$ws.Range("J132").Value = 1370
$ws.Range("K132").Value = 3458.1111
$ws.Range("L132").Value = 4110
$ws.Range("M132").Value = -928.1111000000001
$ws.Range("N132").Value = -9170

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2500
$ws.Range("I3").Value = 2500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2386
$ws.Range("H64").Value = 766.2
$ws.Range("I64").Value = 713.6667
$ws.Range("J64").Value = 845
$ws.Range("K64").Value = 713.6667
$ws.Range("L64").Value = 845
$ws.Range("M64").Value = -488.6667
$ws.Range("N64").Value = -1295
$ws.Range("H67").Value = 766.2
$ws.Range("I67").Value = 713.6667
$ws.Range("J67").Value = 845
$ws.Range("K67").Value = 713.6667
$ws.Range("L67").Value = 845
$ws.Range("M67").Value = 66.33330000000001
$ws.Range("N67").Value = -2405
$ws.Range("H86").Value = 5239.9375
$ws.Range("I86").Value = 2659.889
$ws.Range("J86").Value = 8557.143
$ws.Range("K86").Value = 2659.889
$ws.Range("L86").Value = 8557.143
$ws.Range("M86").Value = -1536.889
$ws.Range("N86").Value = -10803.143
$ws.Range("H89").Value = 5239.9375
$ws.Range("I89").Value = 2659.889
$ws.Range("J89").Value = 8557.143
$ws.Range("K89").Value = 13299.445
$ws.Range("L89").Value = 42785.715
$ws.Range("M89").Value = -7683.445
$ws.Range("N89").Value = -54017.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5617.75
$ws.Range("I31").Value = 3013.4
$ws.Range("J31").Value = 9958.333000000001
$ws.Range("K31").Value = 3013.4
$ws.Range("L31").Value = 9958.333000000001
$ws.Range("M31").Value = -2718.4
$ws.Range("N31").Value = -10548.333
$ws.Range("H34").Value = 5617.75
$ws.Range("I34").Value = 3013.4
$ws.Range("J34").Value = 9958.333000000001
$ws.Range("K34").Value = 3013.4
$ws.Range("L34").Value = 9958.333000000001
$ws.Range("M34").Value = -2811.4
$ws.Range("N34").Value = -10362.333
$ws.Range("H94").Value = 4401.625
$ws.Range("I94").Value = 350
$ws.Range("J94").Value = 4980.4287
$ws.Range("K94").Value = 350
$ws.Range("L94").Value = 4980.4287
$ws.Range("M94").Value = 101
$ws.Range("N94").Value = -5882.4287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 139407.75
$ws.Range("I4").Value = 145037.42
$ws.Range("J4").Value = 100000
$ws.Range("K4").Value = 435112.26
$ws.Range("L4").Value = 300000
$ws.Range("M4").Value = -435000.26
$ws.Range("N4").Value = -300224
$ws.Range("H5").Value = 2292.6155
$ws.Range("I5").Value = 1988
$ws.Range("J5").Value = 2780
$ws.Range("K5").Value = 5964
$ws.Range("L5").Value = 8340
$ws.Range("M5").Value = -5852
$ws.Range("N5").Value = -8564
$ws.Range("H12").Value = 223.1
$ws.Range("I12").Value = 75
$ws.Range("J12").Value = 239.55556
$ws.Range("K12").Value = 225
$ws.Range("L12").Value = 718.66668
$ws.Range("M12").Value = -52
$ws.Range("N12").Value = -1064.66668
$ws.Range("H135").Value = 2292.6155
$ws.Range("I135").Value = 1988
$ws.Range("J135").Value = 2780
$ws.Range("K135").Value = 17892
$ws.Range("L135").Value = 25020
$ws.Range("M135").Value = -15357
$ws.Range("N135").Value = -30090
$ws.Range("H137").Value = 4624.8335
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 4624.8335
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 13874.5005
$ws.Range("N137").Value = -24074.5005
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13587857
$ws.Range("I11").Value = 8183993
$ws.Range("J11").Value = 30571428
$ws.Range("K11").Value = 8183993
$ws.Range("L11").Value = 30571428
$ws.Range("M11").Value = -8183854
$ws.Range("N11").Value = -30571706
$ws.Range("H26").Value = 28740.334
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 28740.334
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 28740.334
$ws.Range("N26").Value = -29300.334
$ws.Range("H50").Value = 28740.334
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 28740.334
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 28740.334
$ws.Range("N50").Value = -29736.334
$ws.Range("H122").Value = 3498.5715
$ws.Range("I122").Value = 3292.3333
$ws.Range("J122").Value = 3653.25
$ws.Range("K122").Value = 9876.999899999999
$ws.Range("L122").Value = 10959.75
$ws.Range("M122").Value = -7426.999899999999
$ws.Range("N122").Value = -15859.75
$ws.Range("H126").Value = 3245.4443
$ws.Range("I126").Value = 3018.3333
$ws.Range("J126").Value = 3699.6667
$ws.Range("K126").Value = 9054.999899999999
$ws.Range("L126").Value = 11099.0001
$ws.Range("M126").Value = -6584.999899999999
$ws.Range("N126").Value = -16039.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1608.3334
$ws.Range("I22").Value = 1216.6666
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1216.6666
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -921.6666
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1608.3334
$ws.Range("I27").Value = 1216.6666
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1216.6666
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1109.6666
$ws.Range("N27").Value = -2214
$ws.Range("H61").Value = 166672420
$ws.Range("I61").Value = 333337000
$ws.Range("J61").Value = 7833.3335
$ws.Range("K61").Value = 333337000
$ws.Range("L61").Value = 7833.3335
$ws.Range("M61").Value = -333336798
$ws.Range("N61").Value = -8237.333500000001
$ws.Range("H113").Value = 166672420
$ws.Range("I113").Value = 333337000
$ws.Range("J113").Value = 7833.3335
$ws.Range("K113").Value = 333337000
$ws.Range("L113").Value = 7833.3335
$ws.Range("M113").Value = -333334830
$ws.Range("N113").Value = -12173.3335
$ws.Range("H122").Value = 2989.2
$ws.Range("I122").Value = 2987.5386
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8962.6158
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6512.6158
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 6838.421
$ws.Range("I132").Value = 6967.5
$ws.Range("J132").Value = 6150
$ws.Range("K132").Value = 20902.5
$ws.Range("L132").Value = 18450
$ws.Range("M132").Value = -18372.5
$ws.Range("N132").Value = -23510
$ws.Range("H136").Value = 868.5625
$ws.Range("I136").Value = 893.2
$ws.Range("J136").Value = 499
$ws.Range("K136").Value = 2679.6
$ws.Range("L136").Value = 1497
$ws.Range("M136").Value = -129.6000000000004
$ws.Range("N136").Value = -6597

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 13152.8
$ws.Range("I51").Value = 21305.6
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 21305.6
$ws.Range("L51").Value = 5000
$ws.Range("M51").Value = -20795.6
$ws.Range("N51").Value = -6020
$ws.Range("H122").Value = 5183.222
$ws.Range("I122").Value = 4069.8
$ws.Range("J122").Value = 6575
$ws.Range("K122").Value = 12209.4
$ws.Range("L122").Value = 19725
$ws.Range("M122").Value = -9759.400000000001
$ws.Range("N122").Value = -24625
